# Generate Report for Handoff
#
# Swaps the "17b0988e..." / "03bcaafe..." rows to their new report
# positions and refreshes the handoff/handback status, dates and error
# detail for the 03bcaafe file (now "Ready for handoff").

function Set-LinkedCellText {
    param($ws, $addr, $text)

    $ws.Range($addr).Value = $text

    $target = $ws.Range($addr).Address()
    $links = @($ws.Hyperlinks)
    foreach ($h in $links) {
        if ($h.Range.Address() -eq $target) {
            $h.TextToDisplay = $text
        }
    }
}

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

Set-LinkedCellText $ov "A2" "17b0988e-966e-4b87-b1bf-6b03ee74d6cc.md"
Set-LinkedCellText $ov "B2" "e2e\17b0988e-966e-4b87-b1bf-6b03ee74d6cc.md"

Set-LinkedCellText $ov "A3" "03bcaafe-f860-45ae-8450-eb80b0940753.md"
Set-LinkedCellText $ov "B3" "e2e\03bcaafe-f860-45ae-8450-eb80b0940753.md"

$ov.Range("E3").Value = "Ready for handoff"
$ov.Range("F3").Value = "Ready for handoff"
$ov.Range("G3").Value = "2016-08-27 22:48:27"

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

Set-LinkedCellText $zh "A2" "17b0988e-966e-4b87-b1bf-6b03ee74d6cc.md"
$zh.Range("G2").Value = "17b0988e-966e-4b87-b1bf-6b03ee74d6cc.46934378ea5b97e4b174f73d16766cb097c1d08d.zh-cn.xlf"
Set-LinkedCellText $zh "I2" "17b0988e-966e-4b87-b1bf-6b03ee74d6cc.md"
$zh.Range("J2").Value = "17b0988e-966e-4b87-b1bf-6b03ee74d6cc.46934378ea5b97e4b174f73d16766cb097c1d08d.zh-cn.xlf"

Set-LinkedCellText $zh "A3" "03bcaafe-f860-45ae-8450-eb80b0940753.md"
$zh.Range("C3").Value = "Ready for handoff"
$zh.Range("G3").Value = "03bcaafe-f860-45ae-8450-eb80b0940753.87f64b8aa6433911bb55d082940dc120e5bacd99.zh-cn.xlf"
$zh.Range("H3").Value = "2016-08-27 22:48:22"
Set-LinkedCellText $zh "I3" "03bcaafe-f860-45ae-8450-eb80b0940753.md"
$zh.Range("J3").Value = "03bcaafe-f860-45ae-8450-eb80b0940753.87f64b8aa6433911bb55d082940dc120e5bacd99.zh-cn.xlf"
$zh.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/fb1c90b214fb8cc698cf41b5b8f3f685ebb05a08/e2e/03bcaafe-f860-45ae-8450-eb80b0940753.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/012ad02ce6e164f06e9b8a32732a3452e8ab60ae/e2e/03bcaafe-f860-45ae-8450-eb80b0940753.md."

$zh.Columns.Item(16).ColumnWidth = 39.1

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

Set-LinkedCellText $de "A2" "17b0988e-966e-4b87-b1bf-6b03ee74d6cc.md"
$de.Range("G2").Value = "17b0988e-966e-4b87-b1bf-6b03ee74d6cc.46934378ea5b97e4b174f73d16766cb097c1d08d.de-de.xlf"
Set-LinkedCellText $de "I2" "17b0988e-966e-4b87-b1bf-6b03ee74d6cc.md"
$de.Range("J2").Value = "17b0988e-966e-4b87-b1bf-6b03ee74d6cc.46934378ea5b97e4b174f73d16766cb097c1d08d.de-de.xlf"

Set-LinkedCellText $de "A3" "03bcaafe-f860-45ae-8450-eb80b0940753.md"
$de.Range("C3").Value = "Ready for handoff"
$de.Range("G3").Value = "03bcaafe-f860-45ae-8450-eb80b0940753.87f64b8aa6433911bb55d082940dc120e5bacd99.de-de.xlf"
$de.Range("H3").Value = "2016-08-27 22:48:27"
Set-LinkedCellText $de "I3" "03bcaafe-f860-45ae-8450-eb80b0940753.md"
$de.Range("J3").Value = "03bcaafe-f860-45ae-8450-eb80b0940753.87f64b8aa6433911bb55d082940dc120e5bacd99.de-de.xlf"
$de.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/fb1c90b214fb8cc698cf41b5b8f3f685ebb05a08/e2e/03bcaafe-f860-45ae-8450-eb80b0940753.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/012ad02ce6e164f06e9b8a32732a3452e8ab60ae/e2e/03bcaafe-f860-45ae-8450-eb80b0940753.md."

$de.Columns.Item(16).ColumnWidth = 39.1
